$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Locals")

# Switch off the pause points (Breakpoint 1 / Breakpoint 2) that were TRUE
$ws.Range("B3").Value = $false
$ws.Range("B4").Value = $false

# Update the active sheet selection / scroll position
$ws.Activate()
$ws.Range("A1").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B4").Select() | Out-Null
